# Applies the weekly update to the Jengibre (ginger) price records.
# Updates Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Precio $/Kg and Kg o Unidades for rows 2-23
# (row 21 is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2023-01-17"
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14429
$ws.Range("P2").Value = 1110

# Row 3
$ws.Range("D3").Value = "2022-07-22"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12500
$ws.Range("P3").Value = 962

# Row 4
$ws.Range("D4").Value = "2022-11-25"
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 1115

# Row 5
$ws.Range("D5").Value = "2021-07-02"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12667
$ws.Range("P5").Value = 974

# Row 6
$ws.Range("D6").Value = "2021-07-12"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("P6").Value = 962

# Row 7
$ws.Range("D7").Value = "2023-02-15"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17429
$ws.Range("N7").Value = "$/caja 15 kilos"
$ws.Range("P7").Value = 1162
$ws.Range("Q7").Value = 15

# Row 8
$ws.Range("D8").Value = "2022-12-19"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14400
$ws.Range("P8").Value = 1108

# Row 9
$ws.Range("D9").Value = "2021-07-29"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("P9").Value = 1346

# Row 10
$ws.Range("D10").Value = "2022-09-28"
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("P10").Value = 1038

# Row 11
$ws.Range("D11").Value = "2022-11-28"
$ws.Range("J11").Value = 900
$ws.Range("M11").Value = 13444
$ws.Range("P11").Value = 1034

# Row 12
$ws.Range("D12").Value = "2022-01-31"
$ws.Range("M12").Value = 12500
$ws.Range("P12").Value = 962

# Row 13
$ws.Range("D13").Value = "2020-11-24"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 23000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 23500
$ws.Range("P13").Value = 1808

# Row 14
$ws.Range("D14").Value = "2021-05-04"
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19500
$ws.Range("P14").Value = 1500

# Row 15
$ws.Range("D15").Value = "2022-01-19"
$ws.Range("J15").Value = 160
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 885

# Row 16
$ws.Range("D16").Value = "2023-02-27"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 16500
$ws.Range("P16").Value = 1269

# Row 17
$ws.Range("D17").Value = "2022-02-24"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 19000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19500
$ws.Range("P17").Value = 1500

# Row 18
$ws.Range("D18").Value = "2021-09-30"
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("P18").Value = 1038

# Row 19
$ws.Range("D19").Value = "2021-07-20"
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 12500
$ws.Range("M19").Value = 12750
$ws.Range("P19").Value = 981

# Row 20
$ws.Range("D20").Value = "2023-03-03"
$ws.Range("J20").Value = 750
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17400
$ws.Range("P20").Value = 1338

# Row 22
$ws.Range("D22").Value = "2022-10-21"
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 10000
$ws.Range("N22").Value = "$/caja 13 kilos"
$ws.Range("P22").Value = 769
$ws.Range("Q22").Value = 13

# Row 23
$ws.Range("D23").Value = "2022-12-23"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12750
$ws.Range("P23").Value = 981

